# Auto-generated edit script: update "F" column ("想去人数") values
# across sheets 展览, 演出, 全部类型 per the target diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(3, 6).Value = 63
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(6, 6).Value = 8714
$ws.Cells.Item(8, 6).Value = 229
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(12, 6).Value = 5301
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(14, 6).Value = 5957
$ws.Cells.Item(15, 6).Value = 1088
$ws.Cells.Item(16, 6).Value = 380
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(18, 6).Value = 538
$ws.Cells.Item(20, 6).Value = 265
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(23, 6).Value = 159
$ws.Cells.Item(26, 6).Value = 0
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(28, 6).Value = 1589
$ws.Cells.Item(29, 6).Value = 46
$ws.Cells.Item(31, 6).Value = 1985
$ws.Cells.Item(32, 6).Value = 76
$ws.Cells.Item(33, 6).Value = 82
$ws.Cells.Item(34, 6).Value = 1023
$ws.Cells.Item(35, 6).Value = 1999
$ws.Cells.Item(36, 6).Value = 0
$ws.Cells.Item(37, 6).Value = 1322
$ws.Cells.Item(39, 6).Value = 0
$ws.Cells.Item(40, 6).Value = 1188
$ws.Cells.Item(41, 6).Value = 615
$ws.Cells.Item(42, 6).Value = 99
$ws.Cells.Item(43, 6).Value = 165
$ws.Cells.Item(44, 6).Value = 0
$ws.Cells.Item(45, 6).Value = 0
$ws.Cells.Item(46, 6).Value = 0
$ws.Cells.Item(47, 6).Value = 1334
$ws.Cells.Item(48, 6).Value = 57
$ws.Cells.Item(49, 6).Value = 1088

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(6, 6).Value = 3
$ws.Cells.Item(7, 6).Value = 30
$ws.Cells.Item(8, 6).Value = 7
$ws.Cells.Item(9, 6).Value = 39
$ws.Cells.Item(10, 6).Value = 189
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(13, 6).Value = 2
$ws.Cells.Item(17, 6).Value = 7
$ws.Cells.Item(21, 6).Value = 2

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 37
$ws.Cells.Item(3, 6).Value = 34
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(8, 6).Value = 8714
$ws.Cells.Item(10, 6).Value = 229
$ws.Cells.Item(11, 6).Value = 7064
$ws.Cells.Item(12, 6).Value = 182
$ws.Cells.Item(19, 6).Value = 5957
$ws.Cells.Item(20, 6).Value = 1088
$ws.Cells.Item(21, 6).Value = 380
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(23, 6).Value = 538
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(26, 6).Value = 142
$ws.Cells.Item(28, 6).Value = 159
$ws.Cells.Item(30, 6).Value = 9816
$ws.Cells.Item(31, 6).Value = 79
$ws.Cells.Item(33, 6).Value = 1589
$ws.Cells.Item(36, 6).Value = 76
$ws.Cells.Item(37, 6).Value = 0
$ws.Cells.Item(39, 6).Value = 1999
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(42, 6).Value = 0
$ws.Cells.Item(44, 6).Value = 615
$ws.Cells.Item(46, 6).Value = 0
$ws.Cells.Item(47, 6).Value = 1094
$ws.Cells.Item(49, 6).Value = 961
$ws.Cells.Item(50, 6).Value = 1334
